$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Berekening oversterfte")

# Add new data row 23 (week 31) that was missing before
$ws.Range("F23").Value = 31
$ws.Range("G23").Value = 2620
$ws.Range("H23").Value = 2944
$ws.Range("I23").Formula = "=G23-H23"

# Extend the totals at row 28 to include the newly added row 23
$ws.Range("G28").Formula = "=SUM(G3:G23)"
$ws.Range("H28").Formula = "=SUM(H3:H23)"

# Update the active selection to reflect where the editor ended up
$ws.Range("I24").Select()
